$wb = $excel.ActiveWorkbook

# Sheets (per workbook.xml order):
#   CypherOutput, Message, CypherOutput_Message, StatOutput, StatOutput_Message
$statOutput = $wb.Worksheets.Item("StatOutput")
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")

# Helper: write a plain numeric-looking string into a cell as TEXT (not a
# Number), matching how the workbook already stores every other value
# (Age, counts, etc.) as shared-string text rather than numeric cells.
# Round-tripping the value through a TEXT() formula and then collapsing it
# back to a static value via copy/paste-values keeps the cell a literal
# string without touching any cell styles / number formats.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace("""", """""")
    $range.Formula = "=TEXT(""" + $escaped + """,""@"")"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Update the result counts on StatOutput row 2
# Columns: A=number_of_files  B=number_of_sample  C=number_of_cases  D=number_of_study
Set-TextValue $statOutput.Range("A2") "25"
Set-TextValue $statOutput.Range("B2") "16"
Set-TextValue $statOutput.Range("C2") "10"
Set-TextValue $statOutput.Range("D2") "2"

# Update the Cypher query text recorded on StatOutput_Message (row 18 = the
# "Cypher" value line of the second repeated info block) — Breed filter
# changed from Akita to Labrador Retriever.
$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Labrador Retriever']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statOutputMessage.Range("A18").Value = $newCypher
